# Watchlist.xlsx: "Test Cases" sheet Runmode column (D) flips from "Y" to
# "N" for rows 3-30 (row 2 stays "Y"), and the window scrolls/selects the
# whole D3:D30 run instead of just D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Flip the Runmode value for rows 3 through 30 to "N".
$ws.Range("D3:D30").Value = "N"

# Update the view: scroll so column C / row 17 anchors the pane, and
# select the whole D3:D30 range (leaves the active cell at D3, the
# top-left cell of the selection).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 17
$ws.Range("D3:D30").Select() | Out-Null
